$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

$values = @{
    "P6" = 1.37
    "AA7" = 6.1
    "AB7" = 16
    "AC7" = 90
    "AD7" = 800
    "AE7" = 7.2
    "AF7" = 11.75
    "AG7" = 10
    "AH7" = 27
    "AI7" = 23
    "AJ7" = 37
    "G7" = 2.67
    "H7" = 3.1
    "I7" = 2.55
    "L7" = 1.39
    "M7" = 2.57
    "N7" = 2.12
    "O7" = 1.57
    "P7" = 1.47
    "Q7" = 2.32
    "R7" = 1.87
    "S7" = 1.75
    "T7" = 7.3
    "U7" = 12.5
    "V7" = 10.25
    "W7" = 30
    "X7" = 25
    "Y7" = 40
    "Z7" = 7.8
    "AA8" = 6.2
    "AB8" = 21
    "AC8" = 150
    "AD8" = 101
    "AE8" = 7.4
    "AF8" = 16
    "AG8" = 13.5
    "AH8" = 50
    "AI8" = 40
    "AJ8" = 65
    "G8" = 2.1
    "H8" = 3.05
    "I8" = 3.45
    "J8" = 1.1
    "K8" = 5.8
    "L8" = 1.52
    "M8" = 2.22
    "N8" = 2.47
    "O8" = 1.42
    "P8" = 1.57
    "Q8" = 2.12
    "R8" = 2.15
    "S8" = 1.55
    "T8" = 5.4
    "U8" = 8.5
    "V8" = 9.75
    "W8" = 19
    "X8" = 22
    "Y8" = 50
    "Z8" = 6.4
    "AB9" = 18
    "AC9" = 120
    "AE9" = 6.3
    "AG9" = 12
    "AH9" = 50
    "AI9" = 40
    "AJ9" = 60
    "G9" = 2.9
    "H9" = 2.37
    "I9" = 3.15
    "J9" = 1.19
    "K9" = 4.2
    "L9" = 1.7
    "M9" = 2.05
    "N9" = 3.05
    "O9" = 1.33
    "P9" = 1.7
    "Q9" = 2.02
    "R9" = 2.2
    "S9" = 1.6
    "T9" = 5.9
    "V9" = 11.25
    "X9" = 35
    "Y9" = 55
    "Z9" = 4.2
    "AB10" = 19
    "AG10" = 12
    "G10" = 2.2
    "I10" = 3.1
    "J10" = 1.08
    "K10" = 8
    "P10" = 1.47
    "U10" = 9.5
    "AB11" = 13.5
    "AC11" = 60
    "AD11" = 450
    "AE11" = 9.25
    "AF11" = 19.5
    "AG11" = 12.5
    "G11" = 1.7
    "I11" = 4.45
    "L11" = 1.32
    "M11" = 3.1
    "N11" = 1.93
    "O11" = 1.7
    "P11" = 1.37
    "Q11" = 2.5
    "R11" = 1.93
    "S11" = 1.78
    "T11" = 5.5
    "U11" = 6.5
    "V11" = 6.9
    "W11" = 11
    "X11" = 11.5
    "Y11" = 22
    "Z11" = 8.75
    "N15" = 2.5
    "O15" = 1.5
    "J20" = 1.06
    "K20" = 10
    "N20" = 2.08
    "O20" = 1.73
    "AE23" = 10
    "G23" = 2.8
    "N23" = 1.65
    "O23" = 2.2
    "P23" = 1.3
    "Q23" = 3.4
    "R23" = 1.57
    "S23" = 2.25
    "Y23" = 26
    "AA32" = 6.4
    "AB32" = 12.5
    "AG32" = 12.5
    "G32" = 1.6
    "H32" = 3.7
    "I32" = 4.65
    "R32" = 1.77
    "S32" = 1.94
    "T32" = 6.3
    "U32" = 6.7
    "W32" = 10
    "X32" = 10.25
    "Y32" = 19
    "P33" = 1.3
    "AA36" = 6.3
    "AE36" = 9.75
    "AF36" = 19.5
    "AI36" = 37
    "G36" = 1.95
    "H36" = 3.2
    "I36" = 3.75
    "L36" = 1.35
    "M36" = 2.72
    "N36" = 2.02
    "O36" = 1.62
    "P36" = 1.4
    "Q36" = 2.5
    "R36" = 1.83
    "S36" = 1.78
    "T36" = 6.4
    "W36" = 17
    "X36" = 17
    "Y36" = 32
    "Z36" = 8.25
    "AA38" = 11.75
    "AB38" = 13.5
    "AC38" = 32
    "AD38" = 150
    "AE38" = 32
    "AF38" = 50
    "AG38" = 18.5
    "AH38" = 110
    "AI38" = 40
    "AJ38" = 30
    "G38" = 1.4
    "H38" = 4.75
    "I38" = 5
    "N38" = 1.29
    "O38" = 3.6
    "R38" = 1.39
    "S38" = 2.55
    "T38" = 14.5
    "U38" = 11.5
    "V38" = 9
    "W38" = 13
    "X38" = 10.25
    "Y38" = 15.5
    "Z38" = 28
    "AB39" = 21
    "J39" = 1.03
    "K39" = 17
    "L39" = 1.17
    "M39" = 5
    "N39" = 1.57
    "O39" = 2.35
    "R39" = 2
    "S39" = 1.75
    "W39" = 8
}

foreach ($cell in $values.Keys) {
    $ws.Range($cell).Value = $values[$cell]
}
